$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) of the last existing data row (60) into the new row 61
$ws.Range("A60:F60").Copy($ws.Range("A61:F61"))

# Now set the actual values for the new row (date stored as Excel serial number)
$ws.Range("A61").Value = 45627
$ws.Range("B61").Value = -0.596
$ws.Range("C61").Value = 0.368
$ws.Range("D61").Value = -0.472
$ws.Range("E61").Value = 0.388
$ws.Range("F61").Value = 1.307
